$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.034.32'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.43%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.821.31'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +3.03%  '

# Row 4
$ws.Range('E4').Value = '  +0.90%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.11'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.75%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.008'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.68%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4311'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.21%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3701'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.97%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07259'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.42%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.128.27'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +22.83%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8703'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.86%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.33'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.12%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.425'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.31%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.648'
$ws.Range('D14').ClearFormats()

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06965'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.12%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.21'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.56%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.015'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.01%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008863'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.35%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.28'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.66%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.083.66'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.74%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.215'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.70%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.03'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.51%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.376.94'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +20.06%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.53'
$ws.Range('D25').ClearFormats()

# Row 26
$ws.Range('E26').Value = '  +1.64%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.42'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.42%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.247'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.29%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.933'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +13.55%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.09'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.83%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08992'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.73%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.182'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +6.28%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7511'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.88%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.434'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.23%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.816'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.43%  '

# Row 36
$ws.Range('E36').Value = '  +0.83%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.126'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.99%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05250'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.26%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01929'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.17%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5126'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.23%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.753'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +9.58%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1655'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.79%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.497'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.20%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.370'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.70%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '107.41'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.22%  '

# Row 46
$ws.Range('E46').Value = '  +3.67%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.009'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.77%  '

# Row 48
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4602'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.68%  '

# Row 49
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.656'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.85%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06234'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.73%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.834'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +6.03%  '

